$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D31").Value = "R로 간트 차트 그리기"
$ws.Range("E31").Value = "https://partrita.github.io/posts/ganttrify/"

$ws.Range("D41").Value = "유튜브가 딥러닝을 도입한 방식 (2)"
$ws.Range("E41").Value = "http://cloudinsight.net/ai/%ec%9c%a0%ed%8a%9c%eb%b8%8c%ea%b0%80-%eb%94%a5%eb%9f%ac%eb%8b%9d%ec%9d%84-%eb%8f%84%ec%9e%85%ed%95%9c-%eb%b0%a9%ec%8b%9d-2/"

$ws.Range("D50").Value = "프로그램 만드는 AI, 알파코드"
$ws.Range("E50").Value = "http://incredible.egloos.com/7535078"

$ws.Range("D52").Value = "[공지] R & Python 과외 안내"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2434&utm_source=rss&utm_medium=rss&utm_campaign=%25ea%25b3%25b5%25ec%25a7%2580-r-python-%25ea%25b3%25bc%25ec%2599%25b8-%25ec%2595%2588%25eb%2582%25b4"
